{"js": "// Kaioken Framework (Estrutura) \u2014 remove the \"Classes para envio de\n// e-mails.\" row from the components/responsibilities table.\n//\n// The commit simply drops the whole table row (icon cell + description\n// cell) that documents the \"envio de e-mails\" component; every other\n// row (and the image each of them anchors) is left untouched. We find\n// the row by its description text rather than a hard-coded index so the\n// script is resilient to any unrelated reshuffling above it.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every cell's body text so we can locate the target row.\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let i = 0; i < rows.items.length; i++) {\n  const cells = rows.items[i].cells;\n  for (let j = 0; j < cells.items.length; j++) {\n    cells.items[j].body.load(\"text\");\n  }\n}\nawait context.sync();\n\nconst targetText = \"Classes para envio de e-mails.\";\nlet targetRow = null;\nfor (let i = 0; i < rows.items.length; i++) {\n  const cells = rows.items[i].cells;\n  for (let j = 0; j < cells.items.length; j++) {\n    if (cells.items[j].body.text.indexOf(targetText) !== -1) {\n      targetRow = rows.items[i];\n      break;\n    }\n  }\n  if (targetRow) break;\n}\n\nif (targetRow) {\n  targetRow.delete();\n  await context.sync();\n}\n", "ps1": "# Kaioken Framework (Estrutura) - remove the \"Classes para envio de\n# e-mails.\" row from the components/responsibilities table.\n#\n# The commit simply drops the whole table row (icon cell + description\n# cell) that documents the \"envio de e-mails\" component; every other\n# row (and the image each of them anchors) is left untouched. We find\n# the row by its description text rather than a hard-coded index so the\n# script is resilient to any unrelated reshuffling above it.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$targetText = \"Classes para envio de e-mails.\"\n$targetIndex = -1\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $row = $t.Rows.Item($i)\n    if ($row.Range.Text -like \"*$targetText*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0) {\n    $t.Rows.Item($targetIndex).Delete()\n}\n"}
